$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.058024686599811
$ws.Cells.Item(2, 4).Value = 1.058013116364222
$ws.Cells.Item(2, 5).Value = 1.071505316094598
$ws.Cells.Item(2, 6).Value = 1.078452525360996
$ws.Cells.Item(2, 9).Value = 1.049543556205747
$ws.Cells.Item(2, 10).Value = 1.063018056746245
$ws.Cells.Item(2, 11).Value = 1.060746512464455
$ws.Cells.Item(2, 12).Value = 1.074202310599899
$ws.Cells.Item(2, 13).Value = 1.081131160339456
$ws.Cells.Item(2, 14).Value = 1.024656386195612
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.05917670877991
$ws.Cells.Item(3, 4).Value = 1.058910360773672
$ws.Cells.Item(3, 5).Value = 1.072646349474519
$ws.Cells.Item(3, 6).Value = 1.079664636935057
$ws.Cells.Item(3, 9).Value = 1.049896084677912
$ws.Cells.Item(3, 10).Value = 1.063821807522098
$ws.Cells.Item(3, 11).Value = 1.061457494957119
$ws.Cells.Item(3, 12).Value = 1.075159108995992
$ws.Cells.Item(3, 13).Value = 1.082160196157026
$ws.Cells.Item(3, 14).Value = 1.024933188014408
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.05992207897531
$ws.Cells.Item(4, 4).Value = 1.059490808882311
$ws.Cells.Item(4, 5).Value = 1.073384971895347
$ws.Cells.Item(4, 6).Value = 1.08044938308939
$ws.Cells.Item(4, 9).Value = 1.05012297768903
$ws.Cells.Item(4, 10).Value = 1.0643412538374
$ws.Cells.Item(4, 11).Value = 1.061916784079721
$ws.Cells.Item(4, 12).Value = 1.075777930463748
$ws.Cells.Item(4, 13).Value = 1.082825903721713
$ws.Cells.Item(4, 14).Value = 1.025111876790004
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.060235418189632
$ws.Cells.Item(5, 4).Value = 1.059734798648679
$ws.Cells.Item(5, 5).Value = 1.07369556075025
$ws.Cells.Item(5, 6).Value = 1.080779393953256
$ws.Cells.Item(5, 9).Value = 1.050218072533775
$ws.Cells.Item(5, 10).Value = 1.064559477829238
$ws.Cells.Item(5, 11).Value = 1.062109686716041
$ws.Cells.Item(5, 12).Value = 1.076038014016435
$ws.Cells.Item(5, 13).Value = 1.083105732662149
$ws.Cells.Item(5, 14).Value = 1.025186896938962
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.060288028411813
$ws.Cells.Item(6, 4).Value = 1.059775763816418
$ws.Cells.Item(6, 5).Value = 1.073747714229304
$ws.Cells.Item(6, 6).Value = 1.080834810329596
$ws.Cells.Item(6, 9).Value = 1.05023402232095
$ws.Cells.Item(6, 10).Value = 1.064596109746018
$ws.Cells.Item(6, 11).Value = 1.062142065227041
$ws.Cells.Item(6, 12).Value = 1.076081679163753
$ws.Cells.Item(6, 13).Value = 1.083152715170818
$ws.Cells.Item(6, 14).Value = 1.025199487253828
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.05992626588437
$ws.Cells.Item(7, 4).Value = 1.059494069205482
$ws.Cells.Item(7, 5).Value = 1.073389121714998
$ws.Cells.Item(7, 6).Value = 1.08045379230257
$ws.Cells.Item(7, 9).Value = 1.050124249494023
$ws.Cells.Item(7, 10).Value = 1.064344170350589
$ws.Cells.Item(7, 11).Value = 1.061919362370811
$ws.Cells.Item(7, 12).Value = 1.075781405983131
$ws.Cells.Item(7, 13).Value = 1.082829642944376
$ws.Cells.Item(7, 14).Value = 1.025112879608169
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.058414032102932
$ws.Cells.Item(8, 4).Value = 1.0583163711907
$ws.Cells.Item(8, 5).Value = 1.071890872520949
$ws.Cells.Item(8, 6).Value = 1.078862075872882
$ws.Cells.Item(8, 9).Value = 1.04966294688547
$ws.Cells.Item(8, 10).Value = 1.063289819747901
$ws.Cells.Item(8, 11).Value = 1.060986950786525
$ws.Cells.Item(8, 12).Value = 1.074525725782257
$ws.Cells.Item(8, 13).Value = 1.081478958800766
$ws.Cells.Item(8, 14).Value = 1.024750019771353
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.055748730186573
$ws.Cells.Item(9, 4).Value = 1.056240116092831
$ws.Cells.Item(9, 5).Value = 1.069253010260551
$ws.Cells.Item(9, 6).Value = 1.076060512610244
$ws.Cells.Item(9, 9).Value = 1.048840744502958
$ws.Cells.Item(9, 10).Value = 1.061427035670071
$ws.Cells.Item(9, 11).Value = 1.059338051503546
$ws.Cells.Item(9, 12).Value = 1.072310796270068
$ws.Cells.Item(9, 13).Value = 1.079097716799012
$ws.Cells.Item(9, 14).Value = 1.024107390924486
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.053971413355758
$ws.Cells.Item(10, 4).Value = 1.054855252868153
$ws.Cells.Item(10, 5).Value = 1.06749589343831
$ws.Cells.Item(10, 6).Value = 1.074194925780419
$ws.Cells.Item(10, 9).Value = 1.048286320491872
$ws.Cells.Item(10, 10).Value = 1.060181856837108
$ws.Cells.Item(10, 11).Value = 1.058234809407174
$ws.Cells.Item(10, 12).Value = 1.070832610522252
$ws.Cells.Item(10, 13).Value = 1.077509388206718
$ws.Cells.Item(10, 14).Value = 1.023676799391414
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.053201691125721
$ws.Cells.Item(11, 4).Value = 1.054255422007438
$ws.Cells.Item(11, 5).Value = 1.066735376094707
$ws.Cells.Item(11, 6).Value = 1.073387595327678
$ws.Cells.Item(11, 9).Value = 1.048044752806856
$ws.Cells.Item(11, 10).Value = 1.059641883075656
$ws.Cells.Item(11, 11).Value = 1.057756143302128
$ws.Cells.Item(11, 12).Value = 1.070192157263854
$ws.Cells.Item(11, 13).Value = 1.076821414372013
$ws.Cells.Item(11, 14).Value = 1.02348983180705
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.05291576049295
$ws.Cells.Item(12, 4).Value = 1.05403259111877
$ws.Cells.Item(12, 5).Value = 1.066452933892593
$ws.Cells.Item(12, 6).Value = 1.073087787794597
$ws.Cells.Item(12, 9).Value = 1.047954798256855
$ws.Cells.Item(12, 10).Value = 1.059441191211836
$ws.Cells.Item(12, 11).Value = 1.057578201220571
$ws.Cells.Item(12, 12).Value = 1.069954204902135
$ws.Cells.Item(12, 13).Value = 1.076565836491641
$ws.Cells.Item(12, 14).Value = 1.023420305638163
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.052977094584761
$ws.Cells.Item(13, 4).Value = 1.0540803903317
$ws.Cells.Item(13, 5).Value = 1.066513516559114
$ws.Cells.Item(13, 6).Value = 1.073152094338941
$ws.Cells.Item(13, 9).Value = 1.047974104027031
$ws.Cells.Item(13, 10).Value = 1.059484245813621
$ws.Cells.Item(13, 11).Value = 1.057616376936691
$ws.Cells.Item(13, 12).Value = 1.070005249201102
$ws.Cells.Item(13, 13).Value = 1.076620660361075
$ws.Cells.Item(13, 14).Value = 1.023435222773547
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.053178056459035
$ws.Cells.Item(14, 4).Value = 1.05423700329073
$ws.Cells.Item(14, 5).Value = 1.066712028365704
$ws.Cells.Item(14, 6).Value = 1.073362811703849
$ws.Cells.Item(14, 9).Value = 1.048037321739471
$ws.Cells.Item(14, 10).Value = 1.059625296308761
$ws.Cells.Item(14, 11).Value = 1.057741437499386
$ws.Cells.Item(14, 12).Value = 1.070172489267378
$ws.Cells.Item(14, 13).Value = 1.076800288920184
$ws.Cells.Item(14, 14).Value = 1.023484086346093
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.053301872800902
$ws.Cells.Item(15, 4).Value = 1.054333494109015
$ws.Cells.Item(15, 5).Value = 1.066834344341139
$ws.Cells.Item(15, 6).Value = 1.073492650988288
$ws.Cells.Item(15, 9).Value = 1.048076242357236
$ws.Cells.Item(15, 10).Value = 1.059712186038159
$ws.Cells.Item(15, 11).Value = 1.057818472313713
$ws.Cells.Item(15, 12).Value = 1.070275523475341
$ws.Cells.Item(15, 13).Value = 1.076910959495103
$ws.Cells.Item(15, 14).Value = 1.023514182453996
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.054022494356942
$ws.Cells.Item(16, 4).Value = 1.054895057911672
$ws.Cells.Item(16, 5).Value = 1.067546373236388
$ws.Cells.Item(16, 6).Value = 1.074248515664661
$ws.Cells.Item(16, 9).Value = 1.048302320936298
$ws.Cells.Item(16, 10).Value = 1.060217676119171
$ws.Cells.Item(16, 11).Value = 1.058266556719689
$ws.Cells.Item(16, 12).Value = 1.070875107008644
$ws.Cells.Item(16, 13).Value = 1.077555042092419
$ws.Cells.Item(16, 14).Value = 1.023689196873396
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.054474484634752
$ws.Cells.Item(17, 4).Value = 1.055247264544731
$ws.Cells.Item(17, 5).Value = 1.067993096573859
$ws.Cells.Item(17, 6).Value = 1.074722777444434
$ws.Cells.Item(17, 9).Value = 1.048443732577502
$ws.Cells.Item(17, 10).Value = 1.06053454103607
$ws.Cells.Item(17, 11).Value = 1.058547372055828
$ws.Cells.Item(17, 12).Value = 1.071251105010269
$ws.Cells.Item(17, 13).Value = 1.077958999158979
$ws.Cells.Item(17, 14).Value = 1.023798839830184
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.054738110260798
$ws.Cells.Item(18, 4).Value = 1.055452683776811
$ws.Cells.Item(18, 5).Value = 1.068253694254401
$ws.Cells.Item(18, 6).Value = 1.074999452810187
$ws.Cells.Item(18, 9).Value = 1.04852607104766
$ws.Cells.Item(18, 10).Value = 1.060719285553683
$ws.Cells.Item(18, 11).Value = 1.058711074794025
$ws.Cells.Item(18, 12).Value = 1.071470380813698
$ws.Cells.Item(18, 13).Value = 1.078194599614513
$ws.Cells.Item(18, 14).Value = 1.023862742684957
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.054827997689978
$ws.Cells.Item(19, 4).Value = 1.055522723570608
$ws.Cells.Item(19, 5).Value = 1.068342556709647
$ws.Cells.Item(19, 6).Value = 1.075093799886049
$ws.Cells.Item(19, 9).Value = 1.048554121833584
$ws.Cells.Item(19, 10).Value = 1.060782265554534
$ws.Cells.Item(19, 11).Value = 1.05876687759743
$ws.Cells.Item(19, 12).Value = 1.071545141872551
$ws.Cells.Item(19, 13).Value = 1.078274929777059
$ws.Cells.Item(19, 14).Value = 1.023884523426702
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.054425991676258
$ws.Cells.Item(20, 4).Value = 1.055209477862177
$ws.Cells.Item(20, 5).Value = 1.067945164135792
$ws.Cells.Item(20, 6).Value = 1.074671888844075
$ws.Cells.Item(20, 9).Value = 1.048428575400999
$ws.Cells.Item(20, 10).Value = 1.060500552444352
$ws.Cells.Item(20, 11).Value = 1.058517252758742
$ws.Cells.Item(20, 12).Value = 1.071210767866164
$ws.Cells.Item(20, 13).Value = 1.077915660555255
$ws.Cells.Item(20, 14).Value = 1.023787081354249
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.053118878803306
$ws.Cells.Item(21, 4).Value = 1.054190885427642
$ws.Cells.Item(21, 5).Value = 1.066653570276141
$ws.Cells.Item(21, 6).Value = 1.073300758744047
$ws.Cells.Item(21, 9).Value = 1.048018711928746
$ws.Cells.Item(21, 10).Value = 1.059583763816659
$ws.Cells.Item(21, 11).Value = 1.057704614266354
$ws.Cells.Item(21, 12).Value = 1.070123242887291
$ws.Cells.Item(21, 13).Value = 1.076747393728023
$ws.Cells.Item(21, 14).Value = 1.023469699398707
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.052296919711259
$ws.Cells.Item(22, 4).Value = 1.053550299443735
$ws.Cells.Item(22, 5).Value = 1.065841769946363
$ws.Cells.Item(22, 6).Value = 1.072439084635962
$ws.Cells.Item(22, 9).Value = 1.047759709284769
$ws.Cells.Item(22, 10).Value = 1.059006638623048
$ws.Cells.Item(22, 11).Value = 1.057192841678915
$ws.Cells.Item(22, 12).Value = 1.069439127968894
$ws.Cells.Item(22, 13).Value = 1.076012662073808
$ws.Cells.Item(22, 14).Value = 1.023269696916439
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.052732668209359
$ws.Cells.Item(23, 4).Value = 1.053889901217308
$ws.Cells.Item(23, 5).Value = 1.066272094884955
$ws.Cells.Item(23, 6).Value = 1.072895835978476
$ws.Cells.Item(23, 9).Value = 1.047897135350524
$ws.Cells.Item(23, 10).Value = 1.059312650623016
$ws.Cells.Item(23, 11).Value = 1.057464221346852
$ws.Cells.Item(23, 12).Value = 1.069801823265386
$ws.Cells.Item(23, 13).Value = 1.076402176182675
$ws.Cells.Item(23, 14).Value = 1.023375764903059
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.054447903595494
$ws.Cells.Item(24, 4).Value = 1.055226552089617
$ws.Cells.Item(24, 5).Value = 1.067966822643589
$ws.Cells.Item(24, 6).Value = 1.07469488306805
$ws.Cells.Item(24, 9).Value = 1.048435424723101
$ws.Cells.Item(24, 10).Value = 1.060515910664872
$ws.Cells.Item(24, 11).Value = 1.058530862656941
$ws.Cells.Item(24, 12).Value = 1.071228994600735
$ws.Cells.Item(24, 13).Value = 1.077935243468907
$ws.Cells.Item(24, 14).Value = 1.023792394657882
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.056437848080899
$ws.Cells.Item(25, 4).Value = 1.056776998478831
$ws.Cells.Item(25, 5).Value = 1.069934700166778
$ws.Cells.Item(25, 6).Value = 1.076784405495674
$ws.Cells.Item(25, 9).Value = 1.049054410756217
$ws.Cells.Item(25, 10).Value = 1.061909192784223
$ws.Cells.Item(25, 11).Value = 1.059765030172507
$ws.Cells.Item(25, 12).Value = 1.072883681284438
$ws.Cells.Item(25, 13).Value = 1.079713468156488
$ws.Cells.Item(25, 14).Value = 1.024273908395151
